# Automatic update of files.
# Rows 16-25 and 27-38 on the active sheet have their data permuted among
# themselves (row 26 is left untouched). For every destination row we need
# the exact content (columns A:AY, including "empty but present" cells)
# that originally lived in a different row.
#
# Plain Value-array read/write loses cells that only carry an empty
# placeholder (no real value), so instead we drive this with real
# Copy/PasteSpecial (same as Excel's clipboard paste), staging each source
# row far below the used range first so that the permutation's cycles never
# read a row that has already been overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"
$stagingFirstRow = 1000

# destination row -> source row (row 26 is unchanged and omitted)
$map = [ordered]@{
    16 = 36
    17 = 25
    18 = 21
    19 = 17
    20 = 24
    21 = 29
    22 = 33
    23 = 31
    24 = 32
    25 = 23
    27 = 35
    28 = 20
    29 = 16
    30 = 22
    31 = 18
    32 = 38
    33 = 19
    34 = 37
    35 = 30
    36 = 27
    37 = 28
    38 = 34
}

# Ordered list of the distinct source rows that need to be staged.
$sourceRows = @($map.Values | Sort-Object -Unique)

# Row in the staging area that holds a copy of a given source row's data.
$stagingRowFor = @{}
$idx = 0
foreach ($sr in $sourceRows) {
    $stagingRowFor[$sr] = $stagingFirstRow + $idx
    $idx = $idx + 1
}

# Step 1: snapshot every needed source row into its staging row, using
# real copy/paste so empty-but-present cells survive the round trip.
foreach ($sr in $sourceRows) {
    $stageRow = $stagingRowFor[$sr]
    $ws.Range("A" + $sr + ":" + $lastCol + $sr).Copy()
    $ws.Range("A" + $stageRow).PasteSpecial(-4104)
}

# Step 2: paste each staged snapshot into its real destination row. Clear
# the destination first: PasteSpecial only overwrites the cells present in
# the copied range, so any cell the destination row had beyond that would
# otherwise survive as a stale leftover.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $stageRow = $stagingRowFor[$srcRow]
    $ws.Range("A" + $destRow + ":" + $lastCol + $destRow).ClearContents()
    $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow).Copy()
    $ws.Range("A" + $destRow).PasteSpecial(-4104)
}

# Step 3: clean up the staging area so it doesn't leak into the sheet.
foreach ($sr in $sourceRows) {
    $stageRow = $stagingRowFor[$sr]
    $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow).ClearContents()
}

$excel.CutCopyMode = 0

Write-Host "Permuted rows 16-25 and 27-38 per mapping."
